$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 417
$ws1.Range("F3").Value = 1043
$ws1.Range("F4").Value = 5844
$ws1.Range("F5").Value = 548
$ws1.Range("F6").Value = 1041
$ws1.Range("F8").Value = 847
$ws1.Range("F11").Value = 613
$ws1.Range("F12").Value = 48
$ws1.Range("F15").Value = 1995
$ws1.Range("F16").Value = 1511
$ws1.Range("F17").Value = 1077
$ws1.Range("F20").Value = 402
$ws1.Range("F21").Value = 636
$ws1.Range("F22").Value = 228
$ws1.Range("F26").Value = 3527
$ws1.Range("F29").Value = 106
$ws1.Range("F30").Value = 160
$ws1.Range("F32").Value = 501
$ws1.Range("F37").Value = 324
$ws1.Range("F38").Value = 831
$ws1.Range("F39").Value = 106
$ws1.Range("F41").Value = 80
$ws1.Range("F42").Value = 87

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 740
$ws2.Range("F6").Value = 403

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1043
$ws4.Range("F5").Value = 5844
$ws4.Range("F6").Value = 548
$ws4.Range("F7").Value = 1041
$ws4.Range("F9").Value = 740
$ws4.Range("F11").Value = 847
$ws4.Range("F13").Value = 403
$ws4.Range("F16").Value = 613
$ws4.Range("F17").Value = 48
$ws4.Range("F21").Value = 1995
$ws4.Range("F22").Value = 1511
$ws4.Range("F23").Value = 1077
$ws4.Range("F26").Value = 402
$ws4.Range("F28").Value = 636
$ws4.Range("F29").Value = 228
$ws4.Range("F31").Value = 3527
$ws4.Range("F34").Value = 106
$ws4.Range("F35").Value = 160
$ws4.Range("F37").Value = 501
$ws4.Range("F41").Value = 324
$ws4.Range("F42").Value = 831
$ws4.Range("F43").Value = 106
$ws4.Range("F45").Value = 80
$ws4.Range("F46").Value = 87
